# Apply the changes described by the diff:
# 1. Rename the worksheet (tab/sheet name) from "alpha4F-HW35.xpc" to "alpha4F"
# 2. Append a new row 16 of averaged-intensity data (Gaussian Quadrature scheme row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet
$ws.Name = "alpha4F"

# 2) Add new row 16 data.
# Copy the formatting of A15 (bold/centered/bordered "index" column style) onto
# A16 first, then overwrite the value - this reuses the existing style (s="1")
# instead of minting a new one.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))
$ws.Cells.Item(16, 1).Value = 14

$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item(16, 3).Value = 1.007487963729129
$ws.Cells.Item(16, 4).Value = 0.9397300125747868
$ws.Cells.Item(16, 5).Value = 1.009179549290269
$ws.Cells.Item(16, 6).Value = 1.007487963729129
$ws.Cells.Item(16, 7).Value = 0.9652064309417071
$ws.Cells.Item(16, 8).Value = 1.031576304979364
$ws.Cells.Item(16, 9).Value = 1.007988454637791
$ws.Cells.Item(16, 10).Value = 0.9397300125747868
$ws.Cells.Item(16, 11).Value = 0.9744547809325279
$ws.Cells.Item(16, 12).Value = 0.9909713723308284
$ws.Cells.Item(16, 13).Value = 0.9935281193588413
